$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green fill used for newly-PASSED test cases (RGB 198,239,206 -> 0xCEEFC6 BGR long)
$greenColor = 13561798

# Rows 2-6 (TC_001, TC_002, TC_003, TC_004, TC_005): mark as PASSED with verified actual
# result text and a green Interior fill.
foreach ($r in 2..6) {
    $ws.Range("H$r").Value = "Agency functionality verified"
    $iCell = $ws.Range("I$r")
    $iCell.Value = "PASSED"
    $iCell.Interior.Color = $greenColor
    $iCell.Interior.PatternColor = $greenColor
}

# TC_006 (row 7): its objective text is replaced with what used to be TC_007's objective
# (the duplicate TC_007 row below is being removed, but this wording is kept).
$ws.Range("C7").Value = "Verify user can edit the agency user created."

# Remove the now-duplicate TC_007 row entirely.
$ws.Rows(8).Delete()
